$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row from diff hunk @@ -3685,25 +3685,25 @@
$ws.Range("H62").Value = 5668.393
$ws.Range("I62").Value = 5518.1577
$ws.Range("J62").Value = 5985.5557
$ws.Range("K62").Value = 5518.1577
$ws.Range("L62").Value = 5985.5557
$ws.Range("M62").Value = -4894.1577
$ws.Range("N62").Value = -7233.5557
# row from diff hunk @@ -3835,25 +3835,25 @@
$ws.Range("H65").Value = 5668.393
$ws.Range("I65").Value = 5518.1577
$ws.Range("J65").Value = 5985.5557
$ws.Range("K65").Value = 27590.7885
$ws.Range("L65").Value = 29927.7785
$ws.Range("M65").Value = -24470.7885
$ws.Range("N65").Value = -36167.7785
# row from diff hunk @@ -4383,22 +4383,22 @@
$ws.Range("H76").Value = 4168
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4168
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 4168
$ws.Range("N76").Value = -4798
# row from diff hunk @@ -4530,22 +4530,22 @@
$ws.Range("H79").Value = 4168
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4168
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 4168
$ws.Range("N79").Value = -6352
# row from diff hunk @@ -6125,22 +6125,22 @@
$ws.Range("H111").Value = 2984.111
$ws.Range("I111").Value = 2032.4
$ws.Range("J111").Value = 4173.75
$ws.Range("K111").Value = 6097.200000000001
$ws.Range("L111").Value = 12521.25
$ws.Range("M111").Value = -3030.200000000001
$ws.Range("N111").Value = -18655.25
# row from diff hunk @@ -6177,25 +6177,25 @@
$ws.Range("H112").Value = 3899.6924
$ws.Range("I112").Value = 2400
$ws.Range("J112").Value = 4024.6667
$ws.Range("K112").Value = 7200
$ws.Range("L112").Value = 12074.0001
$ws.Range("M112").Value = -6092
$ws.Range("N112").Value = -14290.0001
# row from diff hunk @@ -6624,22 +6624,19 @@
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0
# row from diff hunk @@ -7172,22 +7169,22 @@
$ws.Range("H132").Value = 19546.615
$ws.Range("I132").Value = 20883.834
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 62651.50199999999
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -60121.50199999999
$ws.Range("N132").Value = -15560
# row from diff hunk @@ -7423,25 +7420,25 @@
$ws.Range("H137").Value = 38454.312
$ws.Range("I137").Value = 43699.145
$ws.Range("J137").Value = 34375
$ws.Range("K137").Value = 131097.435
$ws.Range("L137").Value = 103125
$ws.Range("M137").Value = -128547.435
$ws.Range("N137").Value = -108225
# row from diff hunk @@ -7475,25 +7472,25 @@
$ws.Range("H138").Value = 1924.72
$ws.Range("I138").Value = 1288.1945
$ws.Range("J138").Value = 3561.5
$ws.Range("K138").Value = 3864.5835
$ws.Range("L138").Value = 10684.5
$ws.Range("M138").Value = 1275.4165
$ws.Range("N138").Value = -20964.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row from diff hunk @@ -11248,22 +11245,22 @@
$ws.Range("H74").Value = 383634
$ws.Range("I74").Value = 857943.1
$ws.Range("J74").Value = 14726.889
$ws.Range("K74").Value = 857943.1
$ws.Range("L74").Value = 14726.889
$ws.Range("M74").Value = -857069.1
$ws.Range("N74").Value = -16474.889
# row from diff hunk @@ -11392,22 +11389,22 @@
$ws.Range("H77").Value = 383634
$ws.Range("I77").Value = 857943.1
$ws.Range("J77").Value = 14726.889
$ws.Range("K77").Value = 4289715.5
$ws.Range("L77").Value = 73634.44499999999
$ws.Range("M77").Value = -4285347.5
$ws.Range("N77").Value = -82370.44499999999
# row from diff hunk @@ -12648,22 +12645,19 @@
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0
# row from diff hunk @@ -14042,22 +14036,22 @@
$ws.Range("H132").Value = 2687.1667
$ws.Range("I132").Value = 907.8333
$ws.Range("J132").Value = 4466.5
$ws.Range("K132").Value = 2723.4999
$ws.Range("L132").Value = 13399.5
$ws.Range("M132").Value = -193.4998999999998
$ws.Range("N132").Value = -18459.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row from diff hunk @@ -15487,25 +15481,25 @@
$ws.Range("H20").Value = 1078.7
$ws.Range("I20").Value = 895
$ws.Range("J20").Value = 1262.4
$ws.Range("K20").Value = 895
$ws.Range("L20").Value = 1262.4
$ws.Range("M20").Value = -648
$ws.Range("N20").Value = -1756.4
# row from diff hunk @@ -19295,22 +19289,22 @@
$ws.Range("H99").Value = 895.3
$ws.Range("I99").Value = 927.375
$ws.Range("J99").Value = 767
$ws.Range("K99").Value = 927.375
$ws.Range("L99").Value = 767
$ws.Range("M99").Value = 570.625
$ws.Range("N99").Value = -3763

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row from diff hunk @@ -21803,22 +21797,22 @@
$ws.Range("H9").Value = 89427.42999999999
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 89427.42999999999
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 89427.42999999999
$ws.Range("N9").Value = -89763.42999999999
# row from diff hunk @@ -22866,22 +22860,22 @@
$ws.Range("H31").Value = 5264194.5
$ws.Range("I31").Value = 6250931
$ws.Range("J31").Value = 1599.3334
$ws.Range("K31").Value = 6250931
$ws.Range("L31").Value = 1599.3334
$ws.Range("M31").Value = -6250636
$ws.Range("N31").Value = -2189.3334
# row from diff hunk @@ -23016,22 +23010,22 @@
$ws.Range("H34").Value = 5264194.5
$ws.Range("I34").Value = 6250931
$ws.Range("J34").Value = 1599.3334
$ws.Range("K34").Value = 6250931
$ws.Range("L34").Value = 1599.3334
$ws.Range("M34").Value = -6250729
$ws.Range("N34").Value = -2003.3334
# row from diff hunk @@ -24373,25 +24367,22 @@
$ws.Range("H62").Value = 5124.25
$ws.Range("I62").Value = 5124.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5124.25
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4500.25
# row from diff hunk @@ -24526,25 +24517,22 @@
$ws.Range("H65").Value = 5124.25
$ws.Range("I65").Value = 5124.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 25621.25
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -22501.25
# row from diff hunk @@ -27815,25 +27803,25 @@
$ws.Range("H132").Value = 64527.312
$ws.Range("I132").Value = 100873.6
$ws.Range("J132").Value = 3950.1667
$ws.Range("K132").Value = 302620.8
$ws.Range("L132").Value = 11850.5001
$ws.Range("M132").Value = -300090.8
$ws.Range("N132").Value = -16910.5001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row from diff hunk @@ -34407,25 +34395,25 @@
$ws.Range("H121").Value = 91929.266
$ws.Range("I121").Value = 119042.336
$ws.Range("J121").Value = 51259.668
$ws.Range("K121").Value = 357127.008
$ws.Range("L121").Value = 153779.004
$ws.Range("M121").Value = -355817.008
$ws.Range("N121").Value = -156399.004
# row from diff hunk @@ -34459,22 +34447,22 @@
$ws.Range("H122").Value = 1524.0834
$ws.Range("I122").Value = 974.5
$ws.Range("J122").Value = 1634
$ws.Range("K122").Value = 8770.5
$ws.Range("L122").Value = 14706
$ws.Range("M122").Value = -6320.5
$ws.Range("N122").Value = -19606
# row from diff hunk @@ -35209,25 +35197,25 @@
$ws.Range("H137").Value = 3462.8572
$ws.Range("I137").Value = 3383.3333
$ws.Range("J137").Value = 3522.5
$ws.Range("K137").Value = 10149.9999
$ws.Range("L137").Value = 10567.5
$ws.Range("M137").Value = -5049.999899999999
$ws.Range("N137").Value = -20767.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row from diff hunk @@ -40140,22 +40128,22 @@
$ws.Range("H97").Value = 1360.7693
$ws.Range("I97").Value = 1192.4117
$ws.Range("J97").Value = 1678.7778
$ws.Range("K97").Value = 1192.4117
$ws.Range("L97").Value = 1678.7778
$ws.Range("M97").Value = -696.4117000000001
$ws.Range("N97").Value = -2670.7778
# row from diff hunk @@ -41543,25 +41531,25 @@
$ws.Range("H126").Value = 2717.4
$ws.Range("I126").Value = 1896.3334
$ws.Range("J126").Value = 3949
$ws.Range("K126").Value = 5689.0002
$ws.Range("L126").Value = 11847
$ws.Range("M126").Value = -3219.0002
$ws.Range("N126").Value = -16787
# row from diff hunk @@ -41837,22 +41825,22 @@
$ws.Range("H132").Value = 4395.2856
$ws.Range("I132").Value = 4409.3335
$ws.Range("J132").Value = 4384.75
$ws.Range("K132").Value = 13228.0005
$ws.Range("L132").Value = 13154.25
$ws.Range("M132").Value = -10698.0005
$ws.Range("N132").Value = -18214.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row from diff hunk @@ -43383,22 +43371,22 @@
$ws.Range("H22").Value = 1569.3846
$ws.Range("I22").Value = 1250.375
$ws.Range("J22").Value = 2079.8
$ws.Range("K22").Value = 1250.375
$ws.Range("L22").Value = 2079.8
$ws.Range("M22").Value = -955.375
$ws.Range("N22").Value = -2669.8
# row from diff hunk @@ -43634,22 +43622,22 @@
$ws.Range("H27").Value = 1569.3846
$ws.Range("I27").Value = 1250.375
$ws.Range("J27").Value = 2079.8
$ws.Range("K27").Value = 1250.375
$ws.Range("L27").Value = 2079.8
$ws.Range("M27").Value = -1143.375
$ws.Range("N27").Value = -2293.8
# row from diff hunk @@ -44084,19 +44072,22 @@
$ws.Range("H36").Value = 75979
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 75979
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 75979
$ws.Range("N36").Value = -77103
# row from diff hunk @@ -48764,25 +48755,25 @@
$ws.Range("H132").Value = 4986.625
$ws.Range("I132").Value = 4966
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 14898
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -12368
$ws.Range("N132").Value = -20057

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row from diff hunk @@ -51389,22 +51380,22 @@
$ws.Range("H45").Value = 18390.125
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 18390.125
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 18390.125
$ws.Range("N45").Value = -19372.125
# row from diff hunk @@ -55286,25 +55277,25 @@
$ws.Range("H126").Value = 458271.47
$ws.Range("I126").Value = 2711.5715
$ws.Range("J126").Value = 1255501.2
$ws.Range("K126").Value = 8134.7145
$ws.Range("L126").Value = 3766503.6
$ws.Range("M126").Value = -5664.7145
$ws.Range("N126").Value = -3771443.6
# row from diff hunk @@ -55580,22 +55571,22 @@
$ws.Range("H132").Value = 31340682
$ws.Range("I132").Value = 62675650
$ws.Range("J132").Value = 5712.5
$ws.Range("K132").Value = 188026950
$ws.Range("L132").Value = 17137.5
$ws.Range("M132").Value = -188024420
$ws.Range("N132").Value = -22197.5
# row from diff hunk @@ -55776,25 +55767,25 @@
$ws.Range("H136").Value = 13733.057
$ws.Range("I136").Value = 15594.578
$ws.Range("J136").Value = 3262
$ws.Range("K136").Value = 46783.734
$ws.Range("L136").Value = 9786
$ws.Range("M136").Value = -44233.734
$ws.Range("N136").Value = -14886
